# Apply updated Q (MVAr) results to res_bus sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 19.07577633059918
$ws.Range("F2").Value = 63.22741186618805
$ws.Range("I2").Value = 15.05842864513397
$ws.Range("J2").Value = 30.88839161396027
$ws.Range("K2").Value = 16.64553451538086
$ws.Range("B3").Value = 16.85542519117371
$ws.Range("F3").Value = 63.63942337036133
$ws.Range("I3").Value = 15.13898193836212
$ws.Range("J3").Value = 32.20508885383606
$ws.Range("K3").Value = 18.5420640707016
$ws.Range("B4").Value = 16.98061855294873
$ws.Range("F4").Value = 63.99045717716217
$ws.Range("I4").Value = 15.04914534091949
$ws.Range("J4").Value = 33.32949483394623
$ws.Range("K4").Value = 18.30558180809021
$ws.Range("B5").Value = 16.87606567061084
$ws.Range("F5").Value = 64.19250726699829
$ws.Range("I5").Value = 14.81450164318085
$ws.Range("J5").Value = 33.91446316242218
$ws.Range("K5").Value = 18.52396762371063
$ws.Range("B6").Value = 16.64085564568813
$ws.Range("F6").Value = 64.28613674640656
$ws.Range("I6").Value = 14.64131593704224
$ws.Range("J6").Value = 34.15962672233582
$ws.Range("K6").Value = 18.87422263622284
$ws.Range("B7").Value = 14.85030156548237
$ws.Range("F7").Value = 64.2725909948349
$ws.Range("I7").Value = 14.64131593704224
$ws.Range("J7").Value = 34.06191146373749
$ws.Range("K7").Value = 20.59826564788818
$ws.Range("B8").Value = 14.38553132440938
$ws.Range("F8").Value = 64.25316953659058
$ws.Range("I8").Value = 14.60904049873352
$ws.Range("J8").Value = 33.97127258777618
$ws.Range("K8").Value = 21.09165012836456
$ws.Range("B9").Value = 15.30500489911719
$ws.Range("F9").Value = 63.98899555206299
$ws.Range("I9").Value = 14.55940771102905
$ws.Range("J9").Value = 33.11190140247345
$ws.Range("K9").Value = 20.47143340110779
$ws.Range("B10").Value = 19.25146876826693
$ws.Range("F10").Value = 63.29695796966553
$ws.Range("I10").Value = 14.72215282917023
$ws.Range("J10").Value = 31.00481188297272
$ws.Range("K10").Value = 16.81749296188354
$ws.Range("B11").Value = 23.49256765301106
$ws.Range("F11").Value = 62.61118912696838
$ws.Range("I11").Value = 14.83576107025146
$ws.Range("J11").Value = 28.88713479042053
$ws.Range("K11").Value = 12.57626569271088
$ws.Range("B12").Value = 24.85759347726707
$ws.Range("F12").Value = 62.28311800956726
$ws.Range("I12").Value = 14.78588438034058
$ws.Range("J12").Value = 27.98708748817444
$ws.Range("K12").Value = 11.07244575023651
$ws.Range("B13").Value = 25.03301543873749
$ws.Range("F13").Value = 62.02623856067657
$ws.Range("I13").Value = 14.73027169704437
$ws.Range("J13").Value = 28.21517086029053
$ws.Range("K13").Value = 9.888022661209106
$ws.Range("B14").Value = 24.41382299921679
$ws.Range("F14").Value = 61.98859691619873
$ws.Range("I14").Value = 14.65024733543396
$ws.Range("J14").Value = 29.01511704921722
$ws.Range("K14").Value = 9.684293508529663
$ws.Range("B15").Value = 23.88713060604277
$ws.Range("F15").Value = 62.06834614276886
$ws.Range("I15").Value = 14.53139305114746
$ws.Range("J15").Value = 29.60489523410797
$ws.Range("K15").Value = 9.995007991790771
$ws.Range("B16").Value = 23.64691552854129
$ws.Range("F16").Value = 62.05498278141022
$ws.Range("I16").Value = 14.6475578546524
$ws.Range("J16").Value = 29.73841512203217
$ws.Range("K16").Value = 9.983206152915955
$ws.Range("B17").Value = 23.42173252939392
$ws.Range("F17").Value = 62.03431844711304
$ws.Range("I17").Value = 14.79836642742157
$ws.Range("J17").Value = 29.82449352741241
$ws.Range("K17").Value = 9.953213214874268
$ws.Range("B18").Value = 22.59416152808262
$ws.Range("F18").Value = 62.25078940391541
$ws.Range("I18").Value = 14.78380405902863
$ws.Range("J18").Value = 30.40887117385864
$ws.Range("K18").Value = 10.92545413970947
$ws.Range("B19").Value = 22.26419300870475
$ws.Range("F19").Value = 62.37529492378235
$ws.Range("I19").Value = 14.85128676891327
$ws.Range("J19").Value = 30.49361562728882
$ws.Range("K19").Value = 11.51665568351746
$ws.Range("B20").Value = 23.42716732119152
$ws.Range("F20").Value = 62.37637031078339
$ws.Range("I20").Value = 14.83804428577423
$ws.Range("J20").Value = 29.29988694190979
$ws.Range("K20").Value = 11.51592206954956
$ws.Range("B21").Value = 25.84974314591454
$ws.Range("F21").Value = 62.02113664150238
$ws.Range("I21").Value = 14.91435158252716
$ws.Range("J21").Value = 27.13181865215302
$ws.Range("K21").Value = 9.942639827728271
$ws.Range("B22").Value = 27.18555382835257
$ws.Range("F22").Value = 61.65805673599243
$ws.Range("I22").Value = 15.08648085594177
$ws.Range("J22").Value = 25.94317972660065
$ws.Range("K22").Value = 8.37472128868103
$ws.Range("B23").Value = 28.07945376405041
$ws.Range("F23").Value = 61.36943030357361
$ws.Range("I23").Value = 15.27212595939636
$ws.Range("J23").Value = 25.08505177497864
$ws.Range("K23").Value = 7.149543881416321
$ws.Range("B24").Value = 28.35238135990949
$ws.Range("F24").Value = 61.25425064563751
$ws.Range("I24").Value = 15.40271854400635
$ws.Range("J24").Value = 24.76133465766907
$ws.Range("K24").Value = 6.685088038444519
$ws.Range("B25").Value = 26.18447377310076
$ws.Range("F25").Value = 61.76413881778717
$ws.Range("I25").Value = 15.49891638755798
$ws.Range("J25").Value = 26.41545653343201
$ws.Range("K25").Value = 9.031234264373779
